$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K:K").Delete()
